{"js": "// 1) \"Microsoft Visual C# (con gestor...\" -> insert \" y Visual Basic\" right\n//    after \"Microsoft Visual C#\" (before the space + \"(con gestor...\").\nconst csharpResults = context.document.body.search(\"Microsoft Visual C#\", { matchCase: true });\ncsharpResults.load(\"items\");\nawait context.sync();\n\nif (csharpResults.items.length > 0) {\n  csharpResults.items[0].insertText(\"Microsoft Visual C# y Visual Basic\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Price correction: \"Q 10000\" + \"0.00\" ( => \"Q 100000.00\" ) becomes\n//    \"Q 15,000\" + \"0.00\" ( => \"Q 15,000.00\" ).\nconst priceResults = context.document.body.search(\"Q 10000\", { matchCase: true });\npriceResults.load(\"items\");\nawait context.sync();\n\nif (priceResults.items.length > 0) {\n  priceResults.items[0].insertText(\"Q 15,000\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"...Microsoft Visual C# (con gestor...\" -> insert \" y Visual Basic\" right\n#    after \"Microsoft Visual C#\" (before the space + \"(con gestor...\").\n$find1 = $d.Content.Find\n$find1.Text = \"Microsoft Visual C#\"\n$find1.Replacement.Text = \"Microsoft Visual C# y Visual Basic\"\n$find1.Execute(\"Microsoft Visual C#\", $false, $true, $false, $false, $false, $true, 1, $false, \"Microsoft Visual C# y Visual Basic\", 2)\n\n# 2) Price correction: \"Q 10000\" + \"0.00\" ( => \"Q 100000.00\" ) becomes\n#    \"Q 15,000\" + \"0.00\" ( => \"Q 15,000.00\" ).\n$find2 = $d.Content.Find\n$find2.Text = \"Q 10000\"\n$find2.Replacement.Text = \"Q 15,000\"\n$find2.Execute(\"Q 10000\", $false, $true, $false, $false, $false, $true, 1, $false, \"Q 15,000\", 2)\n"}
